# Add data for 2024-06-27
#
# A day's worth of violent-crime incident data is folded into the
# 2024 (and, for a handful of late-arriving/reclassified records, the
# 2022/2023) year-to-date columns across the "Citywide Totals" summary,
# the "By Neighborhood" rollup, and every affected per-neighborhood
# sheet. Values below are the new totals after incorporating that day.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 3779
$ws.Range("K3").Value = 3824
$ws.Range("I4").Value = 1796
$ws.Range("K4").Value = 780
$ws.Range("K5").Value = 266
$ws.Range("K6").Value = 4361
$ws.Range("I7").Value = 26250
$ws.Range("K7").Value = 13010

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K2").Value = 41
$ws.Range("K4").Value = 8
$ws.Range("K7").Value = 178

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 252
$ws.Range("K3").Value = 267
$ws.Range("J4").Value = 97
$ws.Range("K6").Value = 288
$ws.Range("J7").Value = 1852
$ws.Range("K7").Value = 877

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K2").Value = 101
$ws.Range("K3").Value = 91
$ws.Range("K7").Value = 278

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 145
$ws.Range("K3").Value = 204
$ws.Range("K5").Value = 11
$ws.Range("K6").Value = 154
$ws.Range("K7").Value = 538

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K3").Value = 76
$ws.Range("K4").Value = 12
$ws.Range("K7").Value = 216

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 120
$ws.Range("K5").Value = 17
$ws.Range("K6").Value = 130
$ws.Range("K7").Value = 438

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K6").Value = 122
$ws.Range("K7").Value = 303

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K3").Value = 92
$ws.Range("K7").Value = 229

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("K3").Value = 19
$ws.Range("K7").Value = 45

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 107
$ws.Range("K4").Value = 44
$ws.Range("K7").Value = 376
$ws.Range("J8").Value = 1852
$ws.Range("K8").Value = 877
$ws.Range("K10").Value = 73
$ws.Range("K14").Value = 70
$ws.Range("K15").Value = 131
$ws.Range("K19").Value = 402
$ws.Range("K20").Value = 293
$ws.Range("K24").Value = 40
$ws.Range("J25").Value = 153
$ws.Range("K27").Value = 129
$ws.Range("K29").Value = 687
$ws.Range("K30").Value = 45
$ws.Range("K31").Value = 142
$ws.Range("K33").Value = 538
$ws.Range("K34").Value = 64
$ws.Range("K37").Value = 438
$ws.Range("I40").Value = 50
$ws.Range("K42").Value = 459
$ws.Range("K44").Value = 121
$ws.Range("K45").Value = 16
$ws.Range("K47").Value = 75
$ws.Range("K49").Value = 72
$ws.Range("K52").Value = 357
$ws.Range("K53").Value = 178
$ws.Range("K55").Value = 147
$ws.Range("K56").Value = 14
$ws.Range("K57").Value = 41
$ws.Range("K58").Value = 6
$ws.Range("J63").Value = 103
$ws.Range("K63").Value = 43
$ws.Range("K65").Value = 303
$ws.Range("K67").Value = 512
$ws.Range("K68").Value = 30
$ws.Range("K73").Value = 117
$ws.Range("K75").Value = 44
$ws.Range("K77").Value = 91
$ws.Range("K79").Value = 335
$ws.Range("K83").Value = 278
$ws.Range("K85").Value = 590
$ws.Range("K86").Value = 88
$ws.Range("K87").Value = 17
$ws.Range("K90").Value = 119
$ws.Range("K91").Value = 140
$ws.Range("K94").Value = 158
$ws.Range("K95").Value = 216
$ws.Range("K97").Value = 109
$ws.Range("K98").Value = 67
$ws.Range("K99").Value = 229
$ws.Range("I101").Value = 26250
$ws.Range("K101").Value = 13010

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K6").Value = 51
$ws.Range("K7").Value = 142

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K6").Value = 149
$ws.Range("K7").Value = 512

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K4").Value = 6
$ws.Range("K7").Value = 72

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 195
$ws.Range("K3").Value = 241
$ws.Range("K6").Value = 197
$ws.Range("K7").Value = 687

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K3").Value = 119
$ws.Range("K4").Value = 16
$ws.Range("K7").Value = 402

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K2").Value = 25
$ws.Range("K7").Value = 121

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("K6").Value = 24
$ws.Range("K7").Value = 70

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 122
$ws.Range("K4").Value = 18
$ws.Range("K6").Value = 168
$ws.Range("K7").Value = 459

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("K4").Value = 5
$ws.Range("K6").Value = 32
$ws.Range("K7").Value = 73

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K2").Value = 45
$ws.Range("K6").Value = 56
$ws.Range("K7").Value = 147

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("K2").Value = 14
$ws.Range("K7").Value = 40

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K3").Value = 66
$ws.Range("K7").Value = 140

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K2").Value = 115
$ws.Range("K3").Value = 108
$ws.Range("K4").Value = 23
$ws.Range("K5").Value = 10
$ws.Range("K6").Value = 79
$ws.Range("K7").Value = 335

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K3").Value = 89
$ws.Range("K7").Value = 293

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 140
$ws.Range("K3").Value = 116
$ws.Range("K6").Value = 91
$ws.Range("K7").Value = 376

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("K2").Value = 18
$ws.Range("K7").Value = 64

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K2").Value = 44
$ws.Range("K3").Value = 28
$ws.Range("K6").Value = 67
$ws.Range("K7").Value = 158

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("J2").Value = 68
$ws.Range("J7").Value = 153

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K6").Value = 25
$ws.Range("K7").Value = 75

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K6").Value = 41
$ws.Range("K7").Value = 131

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("K6").Value = 46
$ws.Range("K7").Value = 67

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K2").Value = 35
$ws.Range("K7").Value = 117

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K6").Value = 34
$ws.Range("K7").Value = 107

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K6").Value = 68
$ws.Range("K7").Value = 109

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K3").Value = 33
$ws.Range("K7").Value = 129

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("K4").Value = 33
$ws.Range("K7").Value = 88

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("K3").Value = 14
$ws.Range("K7").Value = 44

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K5").Value = 6
$ws.Range("K7").Value = 119

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("K6").Value = 8
$ws.Range("K7").Value = 30

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("K6").Value = 23
$ws.Range("K7").Value = 41

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 210
$ws.Range("K3").Value = 196
$ws.Range("K6").Value = 135
$ws.Range("K7").Value = 590

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K6").Value = 11
$ws.Range("K7").Value = 91

$ws = $wb.Worksheets.Item("Jackson Park")
$ws.Range("K6").Value = 8
$ws.Range("K7").Value = 16

$ws = $wb.Worksheets.Item("Magnificent Mile")
$ws.Range("K6").Value = 7
$ws.Range("K7").Value = 14

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("I4").Value = 3
$ws.Range("I7").Value = 50

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K6").Value = 142
$ws.Range("K7").Value = 357

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("K3").Value = 11
$ws.Range("K7").Value = 44

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("K3").Value = 5
$ws.Range("K7").Value = 17

$ws = $wb.Worksheets.Item("Millenium Park")
$ws.Range("K6").Value = 6
$ws.Range("K7").Value = 6
